$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage text-looking-like-numbers so they are written
# as literal text (matching the source data, which is all text/inlineStr),
# without leaving a residual NumberFormat on the real target cells.
$scratch = $ws.Range("ZZ9999")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = "26.933.48"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.737.51"
$ws.Range("E3").Value = "  +1.27%  "
$scratch.Value = "1.001"
$scratch.Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  -0.38%  "
$scratch.Value = "311.27"
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  +0.10%  "
$scratch.Value = "1.001"
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -0.10%  "
$scratch.Value = "0.5023"
$scratch.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  +9.40%  "
$scratch.Value = "0.3580"
$scratch.Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  +4.78%  "
$scratch.Value = "42.26"
$scratch.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  +0.89%  "
$scratch.Value = "0.07256"
$scratch.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  +0.19%  "
$scratch.Value = "1.060"
$scratch.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  +2.02%  "
$scratch.Value = "1.002"
$scratch.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -0.25%  "
$scratch.Value = "20.19"
$scratch.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  +2.54%  "
$scratch.Value = "5.946"
$scratch.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "1.737.17"
$ws.Range("E15").Value = "  +0.85%  "
$scratch.Value = "6.817"
$scratch.Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  -0.46%  "
$scratch.Value = "86.58"
$scratch.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  -1.93%  "
$scratch.Value = "0.00001035"
$scratch.Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  -0.27%  "
$scratch.Value = "0.06437"
$scratch.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("E20").Value = "  -0.05%  "
$scratch.Value = "16.47"
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +0.07%  "
$scratch.Value = "5.727"
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "26.962.74"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("E24").Value = "  +3.94%  "
$scratch.Value = "2.046"
$scratch.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -3.81%  "
$scratch.Value = "153.64"
$scratch.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = "  -0.68%  "
$scratch.Value = "19.87"
$scratch.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").Value = "1.938.36"
$ws.Range("E28").Value = "  +1.08%  "
$scratch.Value = "2.211"
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  +3.84%  "
$scratch.Value = "120.11"
$scratch.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +0.28%  "
$scratch.Value = "1.045"
$scratch.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  +2.33%  "
$scratch.Value = "0.09518"
$scratch.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  +4.90%  "
$scratch.Value = "3.581"
$scratch.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -0.34%  "
$scratch.Value = "5.359"
$scratch.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  +0.59%  "
$scratch.Value = "0.02196"
$scratch.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.28%  "
$scratch.Value = "0.05849"
$scratch.Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  +0.33%  "
$scratch.Value = "11.05"
$scratch.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$scratch.Value = "0.2003"
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$scratch.Value = "4.775"
$scratch.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$scratch.Value = "1.422"
$scratch.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  +1.59%  "
$scratch.Value = "0.6046"
$scratch.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  +2.61%  "
$scratch.Value = "1.111"
$scratch.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -1.51%  "
$scratch.Value = "7.602"
$scratch.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +2.12%  "
$scratch.Value = "12.80"
$scratch.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  +0.70%  "
$scratch.Value = "3.597"
$scratch.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  +0.39%  "
$scratch.Value = "0.5648"
$scratch.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  +0.60%  "
$scratch.Value = "120.03"
$scratch.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +1.14%  "
$scratch.Value = "1.850"
$scratch.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -0.42%  "
$scratch.Value = "0.06663"
$scratch.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +0.04%  "
$scratch.Value = "1.099"
$scratch.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  +1.77%  "
$scratch.Value = "1.001"
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  +0.03%  "

# Remove the scratch cell entirely (contents + formatting) so it leaves no trace
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false
